$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Text content changes: "(percent)" labels become "(in percent)" labels
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "(в процентах)"
$ws.Range("C2").Value = "(in percent)"

# ---------------------------------------------------------------------------
# 2. New column L (2023) added to the table
# ---------------------------------------------------------------------------
$ws.Range("L4").Value = 2023
$ws.Range("L5").Value = 4.9000000000000004
$ws.Range("L6").Value = 4.5999999999999996

$ws.Range("L4").Style = $ws.Range("K4").Style
$ws.Range("L5").Style = $ws.Range("K5").Style
$ws.Range("L6").Style = $ws.Range("K6").Style

# ---------------------------------------------------------------------------
# 3. Column widths: columns A:C become a uniform 41-wide block
# ---------------------------------------------------------------------------
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 40.1

# ---------------------------------------------------------------------------
# 4. Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 13.5
$ws.Rows.Item(6).RowHeight = 28.5

# ---------------------------------------------------------------------------
# 5. Extend the title/percent rows with styled (blank) cells out to column K,
#    matching the rest of the table, and give row 3 its separator styling.
# ---------------------------------------------------------------------------
$titleRow = $ws.Range("D1:K1")
$titleRow.Font.Name = "Times New Roman"
$titleRow.Font.Size = 11
$titleRow.VerticalAlignment = -4108
$titleRow.Borders.LineStyle = -4142

$pctRow = $ws.Range("D2:K2")
$pctRow.Font.Name = "Times New Roman"
$pctRow.Font.Size = 11
$pctRow.VerticalAlignment = -4108
$pctRow.Borders.LineStyle = -4142

$sepRow = $ws.Range("A3:K3")
$sepRow.Font.Name = "Times New Roman"
$sepRow.Font.Size = 11
$sepRow.VerticalAlignment = -4108
$sepRow.Borders.LineStyle = -4142

# ---------------------------------------------------------------------------
# 6. Header row (row 4) switches from top-aligned to vertically centered text
# ---------------------------------------------------------------------------
$ws.Range("A4:L4").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 7. Drop the left indent that used to be on the English title cell (C1)
# ---------------------------------------------------------------------------
$ws.Range("C1").IndentLevel = 0

# ---------------------------------------------------------------------------
# 8. Dimension / selection bookkeeping
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
